$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 36 (Idaho) results were overwritten by a failed run (TimeoutException).
# Clear the numeric result columns B:H to blank text values, matching the
# placeholder blanks already present in K36/L36, and reset B36 back to the
# default (unformatted) style since it loses its date number format.
$ws.Range("B36:H36").Value = "'"
$ws.Range("B36:H36").Style = "Normal"

# Pct Includes Hispanic Black flips to FALSE.
$ws.Range("J36").Value = $false

# Status code now reports the timeout error instead of success.
$ws.Range("O36").Value = "An error occurred. ... TimeoutException('', None, None)"
